$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update final score (K) and MACRO_SCORE (N) values for rows 2-6
$ws.Range("K2").Value = 58.7
$ws.Range("N2").Value = 51.53902399942638

$ws.Range("K3").Value = 57.5
$ws.Range("N3").Value = 51.53902399942638

$ws.Range("K4").Value = 50.7
$ws.Range("N4").Value = 51.53902399942638

$ws.Range("K5").Value = 48.5
$ws.Range("N5").Value = 51.53902399942638

$ws.Range("K6").Value = 45.5
$ws.Range("N6").Value = 51.53902399942638
